$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 changes from shared string "R40" to the new string "1"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "1"
